$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row labels from "_old"/"_new" to "_FV2304"/"_FV2310" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")
$labels  = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value2 = $labels[$i] + "_FV2304"
    $ws.Range($newCols[$i] + "1").Value2 = $labels[$i] + "_FV2310"
}

# --- Turn the used range into an Excel table ("Table1") ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U82"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
